# The author reordered the deck so that the "User diagram" slide (originally
# slide 3) now comes before the "Domain model" slide (originally slide 2).
# Moving slide 3 to position 2 pushes the former slide 2 down to position 3,
# swapping the two slides.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$s.MoveTo(2)
